$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 453
$ws1.Range("F5").Value = 5130
$ws1.Range("F7").Value = 42
$ws1.Range("F9").Value = 328
$ws1.Range("F10").Value = 1

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G4").Value = 79.90000000000001

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 453
$ws4.Range("G5").Value = 79.90000000000001
$ws4.Range("F9").Value = 5130
$ws4.Range("F11").Value = 42
$ws4.Range("F14").Value = 328
$ws4.Range("F15").Value = 1
